$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number (days since 1899-12-30).
# All populated rows (2 through 69) currently store 46074 (2026-02-21) and
# need to be bumped by one day to 46075 (2026-02-22).
for ($r = 2; $r -le 69; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2
    if ($current -eq 46074) {
        $cell.Value = 46075
    }
}
